$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new header/legend row above the table ---
$ws.Rows("1:1").Insert()

# --- Fill the new A1 with the legend text (bold title + regular body) ---
$title = "Gene co-expression network cluster assignments"
$body  = "`nGene clusters were calculated by running the Markov Cluster (MCL) algorithm on the co-expression matrix. Cluster values designate network specific gene clusters and are not compared across networks.`n"
$ws.Range("A1").Value = $title + $body
$ws.Range("A1").Characters(1, $title.Length).Font.Bold = $true
$ws.Range("A1").Characters($title.Length + 1, $body.Length).Font.Size = 11

# --- Row 1 formatting: bottom border under the legend row, wrapped text, taller row ---
$ws.Range("A1:D1").Borders.Item(9).LineStyle = 1
$ws.Range("A1").WrapText = $true
$ws.Rows("1:1").RowHeight = 69.5

# --- Merge the legend cell across the table width ---
$ws.Range("A1:D1").Merge()

# --- Freeze the legend row so the table header stays visible while scrolling ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D15").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
